$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "Exemplo"
$ws.Range("C2").Value = "Teste"
$ws.Range("D2").Value = "exemplo@gmail.com"
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = "Masculino"

# Delete rows 3 through 7 (remove the extra user records)
$ws.Range("A3:H7").EntireRow.Delete()

# Adjust column widths (ColumnWidth -> stored xml width has a constant +5/6 offset
# in this engine, so subtract it to land on the exact target width)
$widthOffset = 0.8333333333333334
$ws.Range("B1").EntireColumn.ColumnWidth = 9 - $widthOffset
$ws.Range("C1").EntireColumn.ColumnWidth = 11 - $widthOffset
$ws.Range("D1").EntireColumn.ColumnWidth = 19 - $widthOffset
